$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.287.28'
$ws.Range("E2").Value = '  +2.00%  '
$ws.Range("D3").Value = '2.098.02'
$ws.Range("E3").Value = '  +0.03%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.67%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '342.55'
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("E6").Value = '  -0.68%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5273'
$ws.Range("E7").Value = '  +2.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4389'
$ws.Range("E8").Value = '  +0.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '55.08'
$ws.Range("E9").Value = '  +2.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09360'
$ws.Range("E10").Value = '  +1.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.175'
$ws.Range("E11").Value = '  +0.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.78'
$ws.Range("E12").Value = '  +0.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.568'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.866'
$ws.Range("E14").Value = '  +1.58%  '
$ws.Range("D15").Value = '1.991.40'
$ws.Range("E15").Value = '  -3.68%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '101.00'
$ws.Range("E16").Value = '  -1.29%  '
$ws.Range("E17").Value = '  +0.67%  '
$ws.Range("E18").Value = '  -0.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '21.14'
$ws.Range("E19").Value = '  +0.64%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06716'
$ws.Range("E20").Value = '  +0.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.378'
$ws.Range("E21").Value = '  +2.88%  '
$ws.Range("E22").Value = '  -0.72%  '
$ws.Range("D23").Value = '30.283.51'
$ws.Range("E23").Value = '  +1.79%  '
$ws.Range("E24").Value = '  -1.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.318'
$ws.Range("E25").Value = '  +0.54%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.009'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.82'
$ws.Range("E27").Value = '  -0.26%  '
$ws.Range("E28").Value = '  +0.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.522'
$ws.Range("E29").Value = '  +1.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.47'
$ws.Range("E30").Value = '  +0.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.135'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.682'
$ws.Range("E32").Value = '  +0.65%  '
$ws.Range("E33").Value = '  +0.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.252'
$ws.Range("E34").Value = '  +0.98%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.905'
$ws.Range("E35").Value = '  -1.43%  '
$ws.Range("E36").Value = '  -3.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02621'
$ws.Range("E37").Value = '  +1.70%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06756'
$ws.Range("E38").Value = '  +0.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '12.66'
$ws.Range("E39").Value = '  +1.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.349'
$ws.Range("E40").Value = '  +0.74%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6969'
$ws.Range("E41").Value = '  -0.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2215'
$ws.Range("E42").Value = '  +0.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6780'
$ws.Range("E43").Value = '  -0.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.30'
$ws.Range("E44").Value = '  +0.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.337'
$ws.Range("E45").Value = '  +0.93%  '
$ws.Range("E46").Value = '  -0.65%  '
$ws.Range("E47").Value = '  +8.45%  '
$ws.Range("E48").Value = '  +0.70%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000350'
$ws.Range("E49").Value = '  -2.92%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.211'
$ws.Range("E50").Value = '  +5.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07307'
$ws.Range("E51").Value = '  +3.65%  '
